$d = $word.ActiveDocument

# Replace the team member "Konečný" with "Hromádka" in the member list.
$d.Content.Find.Execute("Konečný", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Hromádka", 2)
